$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 2
    4 = 2
    5 = 2
    6 = 1
    7 = 0
    8 = 0
    9 = 2
    10 = 2
    11 = 3
    12 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 2
    24 = 2
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 0
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 2
    48 = 3
    49 = 3
    50 = 0
    51 = 0
    52 = 1
    53 = 2
    54 = 2
    55 = 2
    56 = 2
    57 = 3
    58 = 0
    59 = 1
    60 = 1
    61 = 1
    62 = 0
    64 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
